# Scheduled runner update: refresh market-board derived profit columns (H-N)
# across the per-job leve sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below mirror the latest data pull; one previously-blank cell (CRP!M57)
# is cleared and two previously-blank HQ-profit cells (CRP!N81, CRP!N84) are filled in
# now that HQ pricing data became available for those rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 913.5454999999999
$ws.Range("I2").Value = 431.25
$ws.Range("K2").Value = 431.25
$ws.Range("M2").Value = -318.25
$ws.Range("H9").Value = 245.58824
$ws.Range("I9").Value = 260.08334
$ws.Range("K9").Value = 260.08334
$ws.Range("M9").Value = -91.08334000000002
$ws.Range("H69").Value = 40027012
$ws.Range("J69").Value = 45462030
$ws.Range("L69").Value = 136386090
$ws.Range("N69").Value = -136387838
$ws.Range("H72").Value = 40027012
$ws.Range("J72").Value = 45462030
$ws.Range("L72").Value = 409158270
$ws.Range("N72").Value = -409167006
$ws.Range("H132").Value = 1376.8914
$ws.Range("I132").Value = 1376.8914
$ws.Range("K132").Value = 4130.674199999999
$ws.Range("M132").Value = -1600.674199999999
$ws.Range("H135").Value = 802.9231
$ws.Range("I135").Value = 735.2727
$ws.Range("J135").Value = 1175
$ws.Range("K135").Value = 6617.454299999999
$ws.Range("L135").Value = 10575
$ws.Range("M135").Value = -4082.454299999999
$ws.Range("N135").Value = -15645
$ws.Range("H137").Value = 2684.8865
$ws.Range("I137").Value = 2375.6086
$ws.Range("J137").Value = 3023.6191
$ws.Range("K137").Value = 7126.825800000001
$ws.Range("L137").Value = 9070.8573
$ws.Range("M137").Value = -4576.825800000001
$ws.Range("N137").Value = -14170.8573

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4921.9873
$ws.Range("I32").Value = 3400.5593
$ws.Range("K32").Value = 3400.5593
$ws.Range("M32").Value = -3113.5593
$ws.Range("H45").Value = 51754.094
$ws.Range("I45").Value = 85289
$ws.Range("K45").Value = 85289
$ws.Range("M45").Value = -84912
$ws.Range("H61").Value = 2046.4706
$ws.Range("I61").Value = 1986.9375
$ws.Range("K61").Value = 1986.9375
$ws.Range("M61").Value = -1774.9375
$ws.Range("H74").Value = 84549.75999999999
$ws.Range("I74").Value = 71153.414
$ws.Range("K74").Value = 71153.414
$ws.Range("M74").Value = -70279.414
$ws.Range("H77").Value = 84549.75999999999
$ws.Range("I77").Value = 71153.414
$ws.Range("K77").Value = 355767.07
$ws.Range("M77").Value = -351399.07
$ws.Range("H132").Value = 2145.1875
$ws.Range("I132").Value = 1702.5416
$ws.Range("J132").Value = 3473.125
$ws.Range("K132").Value = 5107.6248
$ws.Range("L132").Value = 10419.375
$ws.Range("M132").Value = -2577.6248
$ws.Range("N132").Value = -15479.375
$ws.Range("H136").Value = 2046.4706
$ws.Range("I136").Value = 1986.9375
$ws.Range("K136").Value = 5960.8125
$ws.Range("M136").Value = -3410.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4881.476
$ws.Range("I94").Value = 722.3333
$ws.Range("J94").Value = 15279.333
$ws.Range("K94").Value = 722.3333
$ws.Range("L94").Value = 15279.333
$ws.Range("M94").Value = -271.3333
$ws.Range("N94").Value = -16181.333
$ws.Range("H99").Value = 44762.125
$ws.Range("I99").Value = 64707.375
$ws.Range("K99").Value = 64707.375
$ws.Range("M99").Value = -63209.375
$ws.Range("H134").Value = 2516.652
$ws.Range("I134").Value = 791.69696
$ws.Range("J134").Value = 6895.385
$ws.Range("K134").Value = 2375.09088
$ws.Range("L134").Value = 20686.155
$ws.Range("M134").Value = 159.9091200000003
$ws.Range("N134").Value = -25756.155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2948.75
$ws.Range("I16").Value = 2568.25
$ws.Range("J16").Value = 3329.25
$ws.Range("K16").Value = 2568.25
$ws.Range("L16").Value = 3329.25
$ws.Range("M16").Value = -2281.25
$ws.Range("N16").Value = -3903.25
$ws.Range("H22").Value = 602.1667
$ws.Range("I22").Value = 711
$ws.Range("J22").Value = 493.33334
$ws.Range("K22").Value = 711
$ws.Range("L22").Value = 493.33334
$ws.Range("M22").Value = -361
$ws.Range("N22").Value = -1193.33334
$ws.Range("H31").Value = 3787.0312
$ws.Range("J31").Value = 5333.643
$ws.Range("L31").Value = 5333.643
$ws.Range("N31").Value = -5923.643
$ws.Range("H34").Value = 3787.0312
$ws.Range("J34").Value = 5333.643
$ws.Range("L34").Value = 5333.643
$ws.Range("N34").Value = -5737.643
$ws.Range("H50").Value = 2661.2942
$ws.Range("J50").Value = 2661.2942
$ws.Range("L50").Value = 2661.2942
$ws.Range("N50").Value = -3911.2942
$ws.Range("H57").Value = 52000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 52000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 52000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -53120
$ws.Range("H81").Value = 48499.5
$ws.Range("I81").Value = 31999
$ws.Range("J81").Value = 65000
$ws.Range("K81").Value = 31999
$ws.Range("L81").Value = 65000
$ws.Range("M81").Value = -31001
$ws.Range("N81").Value = -66996
$ws.Range("H82").Value = 34296.668
$ws.Range("J82").Value = 48945
$ws.Range("L82").Value = 48945
$ws.Range("N82").Value = -49667
$ws.Range("H84").Value = 48499.5
$ws.Range("I84").Value = 31999
$ws.Range("J84").Value = 65000
$ws.Range("K84").Value = 95997
$ws.Range("L84").Value = 195000
$ws.Range("M84").Value = -91005
$ws.Range("N84").Value = -204984
$ws.Range("H85").Value = 34296.668
$ws.Range("J85").Value = 48945
$ws.Range("L85").Value = 48945
$ws.Range("N85").Value = -51441
$ws.Range("H86").Value = 8477.333000000001
$ws.Range("I86").Value = 7918.1113
$ws.Range("K86").Value = 7918.1113
$ws.Range("M86").Value = -6795.1113
$ws.Range("H89").Value = 8477.333000000001
$ws.Range("I89").Value = 7918.1113
$ws.Range("K89").Value = 39590.5565
$ws.Range("M89").Value = -33974.5565
$ws.Range("H105").Value = 3544
$ws.Range("I105").Value = 2003.3334
$ws.Range("J105").Value = 4699.5
$ws.Range("K105").Value = 2003.3334
$ws.Range("L105").Value = 4699.5
$ws.Range("M105").Value = -256.3334
$ws.Range("N105").Value = -8193.5
$ws.Range("H113").Value = 2948.75
$ws.Range("I113").Value = 2568.25
$ws.Range("J113").Value = 3329.25
$ws.Range("K113").Value = 2568.25
$ws.Range("L113").Value = 3329.25
$ws.Range("M113").Value = -398.25
$ws.Range("N113").Value = -7669.25
$ws.Range("H132").Value = 55067.156
$ws.Range("I132").Value = 57859.832
$ws.Range("K132").Value = 173579.496
$ws.Range("M132").Value = -171049.496
$ws.Range("H134").Value = 23359.408
$ws.Range("I134").Value = 31316.484
$ws.Range("K134").Value = 93949.452
$ws.Range("M134").Value = -91414.452

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5397991.5
$ws.Range("I4").Value = 9366060
$ws.Range("K4").Value = 28098180
$ws.Range("M4").Value = -28098068
$ws.Range("H34").Value = 628.2857
$ws.Range("I34").Value = 250
$ws.Range("J34").Value = 912
$ws.Range("K34").Value = 750
$ws.Range("L34").Value = 2736
$ws.Range("M34").Value = -666
$ws.Range("N34").Value = -2904
$ws.Range("H37").Value = 42149.4
$ws.Range("J37").Value = 42149.4
$ws.Range("L37").Value = 126448.2
$ws.Range("N37").Value = -126672.2
$ws.Range("H46").Value = 167788.16
$ws.Range("I46").Value = 833932
$ws.Range("J46").Value = 1252.1875
$ws.Range("K46").Value = 2501796
$ws.Range("L46").Value = 3756.5625
$ws.Range("M46").Value = -2501705
$ws.Range("N46").Value = -3938.5625
$ws.Range("H56").Value = 10006329
$ws.Range("I56").Value = 10006329
$ws.Range("K56").Value = 10006329
$ws.Range("M56").Value = -10005799

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1526.0769
$ws.Range("I80").Value = 1534.4
$ws.Range("J80").Value = 1498.3334
$ws.Range("K80").Value = 1534.4
$ws.Range("L80").Value = 1498.3334
$ws.Range("M80").Value = -536.4000000000001
$ws.Range("N80").Value = -3494.3334
$ws.Range("H83").Value = 1526.0769
$ws.Range("I83").Value = 1534.4
$ws.Range("J83").Value = 1498.3334
$ws.Range("K83").Value = 7672
$ws.Range("L83").Value = 7491.666999999999
$ws.Range("M83").Value = -2680
$ws.Range("N83").Value = -17475.667
$ws.Range("H97").Value = 800.3226
$ws.Range("I97").Value = 830.64
$ws.Range("J97").Value = 674
$ws.Range("K97").Value = 830.64
$ws.Range("L97").Value = 674
$ws.Range("M97").Value = -334.64
$ws.Range("N97").Value = -1666
$ws.Range("H126").Value = 3652.2307
$ws.Range("I126").Value = 3409
$ws.Range("J126").Value = 4990
$ws.Range("K126").Value = 10227
$ws.Range("L126").Value = 14970
$ws.Range("M126").Value = -7757
$ws.Range("N126").Value = -19910
$ws.Range("H132").Value = 4463.2915
$ws.Range("I132").Value = 4006.25
$ws.Range("K132").Value = 12018.75
$ws.Range("M132").Value = -9488.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 972.76666
$ws.Range("I82").Value = 1667.8889
$ws.Range("J82").Value = 674.8570999999999
$ws.Range("K82").Value = 1667.8889
$ws.Range("L82").Value = 674.8570999999999
$ws.Range("M82").Value = -1306.8889
$ws.Range("N82").Value = -1396.8571
$ws.Range("H85").Value = 972.76666
$ws.Range("I85").Value = 1667.8889
$ws.Range("J85").Value = 674.8570999999999
$ws.Range("K85").Value = 1667.8889
$ws.Range("L85").Value = 674.8570999999999
$ws.Range("M85").Value = -419.8888999999999
$ws.Range("N85").Value = -3170.8571
$ws.Range("H136").Value = 50417.883
$ws.Range("I136").Value = 57487.527
$ws.Range("J136").Value = 8000
$ws.Range("K136").Value = 172462.581
$ws.Range("L136").Value = 24000
$ws.Range("M136").Value = -169912.581
$ws.Range("N136").Value = -29100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3692.9375
$ws.Range("I122").Value = 2885.7273
$ws.Range("K122").Value = 8657.1819
$ws.Range("M122").Value = -6207.1819
$ws.Range("H132").Value = 208908.66
$ws.Range("I132").Value = 5140.4595
$ws.Range("K132").Value = 15421.3785
$ws.Range("M132").Value = -12891.3785
$ws.Range("H136").Value = 2288.7083
$ws.Range("J136").Value = 6498.5
$ws.Range("L136").Value = 19495.5
$ws.Range("N136").Value = -24595.5
